# Case_0_203 (380 kV case) bus voltage-magnitude results — res_bus/vm_pu.
# Refresh the per-bus vm_pu table (rows 2-25) with the re-run values now that
# the 380 kV case has been completed. Column A (bus index), G and H are
# unaffected by this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for columns B:F (buses 0,1,2,3,4), one row per table row (rows 2-25).
$dataBF = @(
    @(1.02,1.067451379478658,1.070097477037803,1.071374243613339,1.07988788374424),
    @(1.02,1.069150547657915,1.071714694440101,1.072887730851513,1.081483075558775),
    @(1.02,1.07024763503335,1.072759058058904,1.073865064770241,1.082513339489722),
    @(1.02,1.070708291400247,1.073197620747663,1.074275468912099,1.082946010637905),
    @(1.02,1.07078560520406,1.073271229072656,1.074344350481441,1.083018631841976),
    @(1.02,1.070253792526133,1.072764920056201,1.073870550431478,1.08251912262331),
    @(1.02,1.068026123358762,1.070644461529055,1.071886152658066,1.080427392738543),
    @(1.02,1.064081831432633,1.0668914768118,1.068373661969892,1.076726220884538),
    @(1.02,1.061438834366485,1.064377708906285,1.066020779956245,1.074247826079461),
    @(1.02,1.060291015499069,1.063286267698655,1.0649991496226,1.073171905976363),
    @(1.02,1.059864140742476,1.062880398928602,1.064619234851172,1.072771833735838),
    @(1.02,1.059955730774659,1.062967480146158,1.06470074778236,1.072857670151978),
    @(1.02,1.060255740671729,1.063252727901081,1.064967754711996,1.073138844679698),
    @(1.02,1.060440516810347,1.06342841723402,1.065132208370635,1.073312028528672),
    @(1.02,1.06151493872634,1.064450080666139,1.066088521778084,1.074319172110873),
    @(1.02,1.06218797850155,1.065090140182921,1.066687628916396,1.074950178835594),
    @(1.02,1.062580225783435,1.065463191107527,1.067036807040591,1.075317969092509),
    @(1.02,1.062713917233339,1.065590343910984,1.067155822180567,1.075443331462884),
    @(1.02,1.062115801475197,1.065021497428497,1.066623378505117,1.074882505320886),
    @(1.02,1.060167409816035,1.063168742326012,1.064889139957385,1.073056057699602),
    @(1.02,1.058939340828691,1.062001179644496,1.063796227755254,1.071905216972483),
    @(1.02,1.059590656553422,1.062620383662864,1.06437584500941,1.072515538969009),
    @(1.02,1.062148416172056,1.065052515015018,1.066652411337513,1.074913084888939),
    @(1.02,1.065103840898896,1.06786373701346,1.069283655083012,1.077684942046455)
)

# Values for columns I:N (buses 7,8,9,10,11,12), one row per table row (rows 2-25).
$dataIN = @(
    @(1.044889200708512,1.072396211224438,1.072798223469824,1.074071586980101,1.082562757274905,1.073919136849941),
    @(1.045273718352527,1.073748229705646,1.074229755469942,1.075399894779469,1.083974217330171,1.075273075352471),
    @(1.045520103889545,1.07462033758161,1.075153484243903,1.076256898880602,1.084885100977442,1.076146421721651),
    @(1.045623107621272,1.074986326054988,1.075541215317997,1.076616594427477,1.085267465713584,1.076512929940641),
    @(1.045640368676892,1.07504773955404,1.075606281921402,1.07667695462911,1.085331633252266,1.076574430653906),
    @(1.045521482492013,1.074625230460864,1.075158667486516,1.07626170745125,1.084890212382729,1.076151321549355),
    @(1.045019653985815,1.072853705217672,1.073282554399747,1.074521018171667,1.083040274980436,1.074377280537121),
    @(1.044116659085397,1.069710630851688,1.069956469618079,1.071434127464825,1.079761427591186,1.071229742645022),
    @(1.043501870772763,1.067600217765132,1.067724889738581,1.069362438726926,1.077562113978871,1.069116332529815),
    @(1.043232578723023,1.066682679933749,1.06675507905504,1.068461975871185,1.076606464625009,1.068197491689467),
    @(1.043132084034106,1.066341295166106,1.066394306366006,1.06812698016658,1.076250981114761,1.067855622116298),
    @(1.043153661715203,1.066414549365351,1.066471718056281,1.068198861770358,1.076327256936279,1.067928980344908),
    @(1.043224281364674,1.066654472668764,1.066725268598931,1.068434295773319,1.076577090778879,1.068169244366932),
    @(1.04326773036939,1.066802221387935,1.066881417277842,1.068579284684509,1.076730953338428,1.06831720290623),
    @(1.043519677451929,1.067661032417341,1.067789177596973,1.069422126683423,1.07762546603936,1.069177233545806),
    @(1.043676888402646,1.06819873867733,1.068357639833858,1.069949898572486,1.078185669575975,1.069715703410359),
    @(1.043768289552184,1.068512016026414,1.068688875294423,1.070257410964056,1.078512106156795,1.070029425649213),
    @(1.043799404674626,1.068618775216617,1.068801760917797,1.070362209456493,1.078623358571446,1.070136336449711),
    @(1.043660051940885,1.068141084962549,1.068296684366106,1.069893307656983,1.078125598294086,1.069657967820691),
    @(1.043203498580908,1.06658383707945,1.066650619343238,1.06836498087892,1.076503535193198,1.068098508466979),
    @(1.042913737678878,1.065601427433621,1.065612531739378,1.067401025284257,1.075480706768446,1.067114703686769),
    @(1.043067603485855,1.066122538867172,1.066163143641838,1.067912328088151,1.07602321337407,1.067636555158327),
    @(1.043667660532001,1.068167137299557,1.06832422856547,1.06991887966932,1.078152742911454,1.069684057155005),
    @(1.044352343790768,1.070525793652627,1.070818795072532,1.072234543838418,1.080611406477866,1.072046063070638)
)

for ($i = 0; $i -lt $dataBF.Length; $i++) {
    $row = 2 + $i

    $rowBF = $dataBF[$i]
    for ($j = 0; $j -lt $rowBF.Length; $j++) {
        # columns B(2) through F(6)
        $ws.Cells.Item($row, 2 + $j).Value = $rowBF[$j]
    }

    $rowIN = $dataIN[$i]
    for ($k = 0; $k -lt $rowIN.Length; $k++) {
        # columns I(9) through N(14)
        $ws.Cells.Item($row, 9 + $k).Value = $rowIN[$k]
    }
}
